# Auto-applies the meteocat daily-summary refresh described by the commit
# message "Update automàtic: dades i banners [2026-02-10 19:20]".
# Re-extracted rows get a later DATA_EXTRACCIO timestamp plus refreshed
# observation figures (precipitation, humidity, pressure, wind gust, temps).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Writing the literal text straight through .Value lets Excel's COM layer
# "smart"-convert look-alike strings (e.g. "93%" -> 0.93 formatted as a
# percentage, picking up a brand-new number-format style in the process).
# Every cell in this sheet is plain text, so instead we build the text via
# a throwaway concatenation formula, then immediately flatten that formula
# down to its literal value with Copy + PasteSpecial(xlPasteValues). That
# keeps the original style index intact and leaves a plain string behind.
function Set-LiteralText {
    param($Range, [string]$Text)
    $escaped = $Text -replace '"', '""'
    $Range.Formula = '="' + $escaped + '"'
    $Range.Copy() | Out-Null
    $Range.PasteSpecial(-4163) | Out-Null
}

$updates = @(
    @{ Cell = 'E2'; Value = '2026-02-10 19:18:33' }
    @{ Cell = 'I2'; Value = '39.3 mm' }
    @{ Cell = 'E3'; Value = '2026-02-10 19:18:35' }
    @{ Cell = 'I3'; Value = '24.1 mm' }
    @{ Cell = 'E4'; Value = '2026-02-10 19:18:38' }
    @{ Cell = 'E5'; Value = '2026-02-10 19:18:40' }
    @{ Cell = 'I5'; Value = '31.8 mm' }
    @{ Cell = 'E6'; Value = '2026-02-10 19:18:42' }
    @{ Cell = 'H6'; Value = '93%' }
    @{ Cell = 'E7'; Value = '2026-02-10 19:18:45' }
    @{ Cell = 'J7'; Value = '1004.6 hPa' }
    @{ Cell = 'O7'; Value = '15.2 °C' }
    @{ Cell = 'E8'; Value = '2026-02-10 19:18:47' }
    @{ Cell = 'J8'; Value = '1004.5 hPa' }
    @{ Cell = 'O8'; Value = '11.9 °C' }
    @{ Cell = 'E9'; Value = '2026-02-10 19:18:50' }
    @{ Cell = 'E10'; Value = '2026-02-10 19:18:52' }
    @{ Cell = 'E11'; Value = '2026-02-10 19:18:55' }
    @{ Cell = 'I11'; Value = '0.8 mm' }
    @{ Cell = 'E12'; Value = '2026-02-10 19:18:58' }
    @{ Cell = 'E13'; Value = '2026-02-10 19:19:00' }
    @{ Cell = 'I13'; Value = '5.4 mm' }
    @{ Cell = 'E14'; Value = '2026-02-10 19:19:03' }
    @{ Cell = 'E15'; Value = '2026-02-10 19:19:06' }
    @{ Cell = 'O15'; Value = '9.0 °C' }
    @{ Cell = 'E16'; Value = '2026-02-10 19:19:08' }
    @{ Cell = 'I16'; Value = '24.9 mm' }
    @{ Cell = 'E17'; Value = '2026-02-10 19:19:10' }
    @{ Cell = 'E18'; Value = '2026-02-10 19:19:12' }
    @{ Cell = 'E19'; Value = '2026-02-10 19:19:15' }
    @{ Cell = 'O19'; Value = '6.5 °C' }
    @{ Cell = 'E20'; Value = '2026-02-10 19:19:17' }
    @{ Cell = 'I20'; Value = '7.6 mm' }
    @{ Cell = 'L20'; Value = '68.4 km/h - 272º 18:52 TU' }
    @{ Cell = 'O20'; Value = '0.4 °C' }
    @{ Cell = 'E21'; Value = '2026-02-10 19:19:20' }
    @{ Cell = 'I21'; Value = '7.9 mm' }
    @{ Cell = 'J21'; Value = '1006.2 hPa' }
    @{ Cell = 'O21'; Value = '7.2 °C' }
    @{ Cell = 'E22'; Value = '2026-02-10 19:19:23' }
    @{ Cell = 'I22'; Value = '9.1 mm' }
    @{ Cell = 'L22'; Value = '86.0 km/h - 310º 18:53 TU' }
    @{ Cell = 'O22'; Value = '-0.5 °C' }
    @{ Cell = 'E23'; Value = '2026-02-10 19:19:25' }
    @{ Cell = 'I23'; Value = '24.8 mm' }
    @{ Cell = 'E24'; Value = '2026-02-10 19:19:27' }
    @{ Cell = 'J24'; Value = '1006.1 hPa' }
    @{ Cell = 'E25'; Value = '2026-02-10 19:19:30' }
    @{ Cell = 'I25'; Value = '17.6 mm' }
    @{ Cell = 'E26'; Value = '2026-02-10 19:19:32' }
    @{ Cell = 'J26'; Value = '1003.4 hPa' }
    @{ Cell = 'O26'; Value = '6.1 °C' }
    @{ Cell = 'E27'; Value = '2026-02-10 19:19:34' }
    @{ Cell = 'I27'; Value = '6.8 mm' }
    @{ Cell = 'E28'; Value = '2026-02-10 19:19:37' }
    @{ Cell = 'E29'; Value = '2026-02-10 19:19:39' }
    @{ Cell = 'O29'; Value = '10.6 °C' }
    @{ Cell = 'E30'; Value = '2026-02-10 19:19:42' }
    @{ Cell = 'L30'; Value = '19.8 km/h - 306º 18:52 TU' }
    @{ Cell = 'E31'; Value = '2026-02-10 19:19:45' }
    @{ Cell = 'E32'; Value = '2026-02-10 19:19:47' }
    @{ Cell = 'O32'; Value = '10.4 °C' }
    @{ Cell = 'E33'; Value = '2026-02-10 19:19:50' }
    @{ Cell = 'I33'; Value = '10.2 mm' }
    @{ Cell = 'E34'; Value = '2026-02-10 19:19:53' }
    @{ Cell = 'I34'; Value = '10.6 mm' }
    @{ Cell = 'E35'; Value = '2026-02-10 19:19:56' }
    @{ Cell = 'H35'; Value = '69%' }
    @{ Cell = 'J35'; Value = '1005.0 hPa' }
    @{ Cell = 'O35'; Value = '13.0 °C' }
    @{ Cell = 'E36'; Value = '2026-02-10 19:19:58' }
    @{ Cell = 'H36'; Value = '93%' }
    @{ Cell = 'O36'; Value = '10.0 °C' }
    @{ Cell = 'E37'; Value = '2026-02-10 19:20:01' }
    @{ Cell = 'I37'; Value = '0.1 mm' }
    @{ Cell = 'J37'; Value = '1005.7 hPa' }
    @{ Cell = 'O37'; Value = '6.6 °C' }
    @{ Cell = 'E38'; Value = '2026-02-10 19:20:06' }
    @{ Cell = 'O38'; Value = '10.8 °C' }
    @{ Cell = 'E39'; Value = '2026-02-10 19:20:09' }
    @{ Cell = 'I39'; Value = '8.0 mm' }
    @{ Cell = 'E40'; Value = '2026-02-10 19:20:15' }
    @{ Cell = 'H40'; Value = '90%' }
    @{ Cell = 'I40'; Value = '12.8 mm' }
    @{ Cell = 'O40'; Value = '7.5 °C' }
    @{ Cell = 'E41'; Value = '2026-02-10 19:20:18' }
    @{ Cell = 'J41'; Value = '1004.7 hPa' }
    @{ Cell = 'E42'; Value = '2026-02-10 19:20:20' }
    @{ Cell = 'E43'; Value = '2026-02-10 19:20:23' }
    @{ Cell = 'H43'; Value = '87%' }
    @{ Cell = 'O43'; Value = '9.5 °C' }
    @{ Cell = 'E44'; Value = '2026-02-10 19:20:25' }
    @{ Cell = 'I44'; Value = '26.1 mm' }
    @{ Cell = 'E45'; Value = '2026-02-10 19:20:28' }
    @{ Cell = 'I45'; Value = '33.2 mm' }
    @{ Cell = 'K45'; Value = '4.1 MJ/m2' }
    @{ Cell = 'E46'; Value = '2026-02-10 19:20:30' }
    @{ Cell = 'H46'; Value = '83%' }
    @{ Cell = 'J46'; Value = '1006.0 hPa' }
    @{ Cell = 'L46'; Value = '33.8 km/h - 279º 18:59 TU' }
    @{ Cell = 'O46'; Value = '14.2 °C' }
)

foreach ($u in $updates) {
    Set-LiteralText $ws.Range($u.Cell) $u.Value
}

$excel.CutCopyMode = $false
